$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text formatting (Price/Volume columns are
# stored as plain text in the source data, even when they look numeric).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '39.442.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.71%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.160.33'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.02%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.62'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.53%  '

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.04%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '64.13'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.63%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0857'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.28%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.00'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.68%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.479.82'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.03%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.21'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.60%  '

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.62%  '

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.97%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.167.44'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.58%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '39.424.24'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.67%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.80'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.21%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.62%  '

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.22%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '231.11'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.44%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.00%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.33'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.12%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.39%  '

$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.65'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.32%  '

$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '172.33'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.41%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.139'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.88'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.80%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.78%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.66'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +6.14%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.57%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.62'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.89%  '

$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'THORChain'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.12'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +9.66%  '

$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.74'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.63%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0617'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.40%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.41'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.62%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.56'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.54%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.02%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '103.84'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.29%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.88%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.75'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.96%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.540.02'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.45%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.92%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.95%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0927'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.58%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.82%  '

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.29%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.75'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.21%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.363.19'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.17%  '

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.13%  '
